$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document has four "Main Success Scenario" sections (Customer Price
# Update, Request Delivery, Update Transport Costs, Discontinue Routes) that
# each end with:
#   "System adds <X> event to log file"
#   "System updates business figures"
# and a following "Exception Scenarios" section with a single line:
#   "1" "a." <tab> "N/A"
#
# The edit adds a new "System verifies ..." step before the "adds ... event
# to log file" step (splitting it into two paragraphs) and replaces the
# placeholder exception-scenario line with real text (renumbering it too).
#
# We process the four blocks from the bottom of the document upward so that
# inserting new paragraphs never invalidates the paragraph index of a block
# we still have to process.
# ---------------------------------------------------------------------------

$wdReplaceAll = 2
$tab = [char]9

function Split-LogParagraph($paraIndex, $oldText, $verifyText) {
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    $replacement = $verifyText + "^p" + $oldText
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, $wdReplaceAll) | Out-Null
}

# --- Block 4: Discontinue Routes (numId=5) ---------------------------------
# Para 102: "System adds discontinue route event to log file"
Split-LogParagraph 102 "System adds discontinue route event to log file" "System verifies that fields are valid"

# The split above inserted one extra paragraph, so the exception line that
# used to be Paragraph 105 is now Paragraph 106.
$para = $d.Paragraphs(106)
$para.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "8", $wdReplaceAll) | Out-Null

$para = $d.Paragraphs(106)
$bmStart = $para.Range.Start + 1
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$para = $d.Paragraphs(106)
$para.Range.Find.Execute("a." + $tab + "N/A", $false, $false, $false, $false, $false, $true, 1, $false, "a. Clerk / Manager enters invalid data in at least one of the fields, go to step 2", $wdReplaceAll) | Out-Null

# --- Block 3: Update Transport Costs (numId=4) ------------------------------
Split-LogParagraph 82 "System adds transport price update event to log file" "System verifies if fields are correct"

# The split above inserted one extra paragraph, so the exception line that
# used to be Paragraph 85 is now Paragraph 86.
$para = $d.Paragraphs(86)
$para.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "13", $wdReplaceAll) | Out-Null

$para = $d.Paragraphs(86)
$para.Range.Find.Execute("a." + $tab + "N/A", $false, $false, $false, $false, $false, $true, 1, $false, "a. Clerk / Manager enters invalid data in at least one of the fields, go to step 2", $wdReplaceAll) | Out-Null

# --- Block 2: Request Delivery (numId=2) ------------------------------------
Split-LogParagraph 56 "System adds delivery event to log file" "System verifies if fields are correct"

# The split above inserted one extra paragraph, so the exception line that
# used to be Paragraph 59 is now Paragraph 60.
$para = $d.Paragraphs(60)
$para.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "9", $wdReplaceAll) | Out-Null

$para = $d.Paragraphs(60)
$para.Range.Find.Execute("N/A", $true, $false, $false, $false, $false, $true, 1, $false, "Clerk / Manager enters invalid data in at least one of the fields, go to step 2", $wdReplaceAll) | Out-Null

# --- Block 1: Customer Price Update (numId=1) -------------------------------
Split-LogParagraph 35 "System adds customer price update event to log file" "System verifies if the fields are correct"

# The split above inserted one extra paragraph, so the exception line that
# used to be Paragraph 38 is now Paragraph 39.
$para = $d.Paragraphs(39)
$para.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "9", $wdReplaceAll) | Out-Null

$para = $d.Paragraphs(39)
$para.Range.Find.Execute("a." + $tab + "N/A", $false, $false, $false, $false, $false, $true, 1, $false, "a. Clerk / Manager enters invalid data in at least one of the fields, go to step 2", $wdReplaceAll) | Out-Null

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
